$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values look numeric ("D.DD" style) as Text so
# Excel keeps them as strings instead of silently converting to numbers.
# (Grouped into contiguous ranges so only the touched rows get re-styled and
# only a single extra style slot is introduced.)
$ws.Range("D4:D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19:D34").NumberFormat = "@"
$ws.Range("D36:D38").NumberFormat = "@"
$ws.Range("D40:D41").NumberFormat = "@"
$ws.Range("D43:D45").NumberFormat = "@"
$ws.Range("D48:D51").NumberFormat = "@"

# Apply the updated values (prices in column D, 1h volume % in column E,
# plus the Cosmos / EthereumClassic row swap in B28:C29).
$ws.Range("D2").Value = '25.986.64'
$ws.Range("E2").Value = '  -1.98%  '
$ws.Range("D3").Value = '1.638.44'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("D4").Value = '1.021'
$ws.Range("E4").Value = '  +1.41%  '
$ws.Range("D5").Value = '216.39'
$ws.Range("E5").Value = '  -1.46%  '
$ws.Range("D6").Value = '0.5065'
$ws.Range("E6").Value = '  -1.73%  '
$ws.Range("D7").Value = '1.020'
$ws.Range("E7").Value = '  +1.35%  '
$ws.Range("D8").Value = '0.2582'
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").Value = '0.06422'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("D10").Value = '19.48'
$ws.Range("E10").Value = '  -2.82%  '
$ws.Range("D11").Value = '0.07783'
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").Value = '1.645.46'
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("D13").Value = '4.261'
$ws.Range("E13").Value = '  -2.05%  '
$ws.Range("D14").Value = '1.865.31'
$ws.Range("E14").Value = '  -1.90%  '
$ws.Range("D15").Value = '0.5458'
$ws.Range("E15").Value = '  -2.27%  '
$ws.Range("D16").Value = '0.0₅7969'
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("E17").Value = '  -1.88%  '
$ws.Range("D18").Value = '25.997.29'
$ws.Range("E18").Value = '  -2.06%  '
$ws.Range("D19").Value = '1.020'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").Value = '204.33'
$ws.Range("E20").Value = '  -3.35%  '
$ws.Range("D21").Value = '4.318'
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("D22").Value = '10.02'
$ws.Range("E22").Value = '  -1.18%  '
$ws.Range("D23").Value = '5.969'
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("D24").Value = '1.021'
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").Value = '1.984'
$ws.Range("E25").Value = '  +14.08%  '
$ws.Range("D26").Value = '141.94'
$ws.Range("E26").Value = '  -1.79%  '
$ws.Range("D27").Value = '0.1156'
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '6.823'
$ws.Range("E28").Value = '  -2.86%  '
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '15.72'
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("D30").Value = '1.245'
$ws.Range("E30").Value = '  -1.48%  '
$ws.Range("D31").Value = '0.04995'
$ws.Range("E31").Value = '  -4.44%  '
$ws.Range("D32").Value = '3.266'
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("D33").Value = '3.207'
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").Value = '1.538'
$ws.Range("E34").Value = '  -2.84%  '
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("D36").Value = '2.636'
$ws.Range("E36").Value = '  -4.96%  '
$ws.Range("D37").Value = '0.8925'
$ws.Range("E37").Value = '  -3.68%  '
$ws.Range("D38").Value = '0.5664'
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").Value = '1.116.58'
$ws.Range("E39").Value = '  -4.15%  '
$ws.Range("D40").Value = '0.01569'
$ws.Range("E40").Value = '  -2.03%  '
$ws.Range("D41").Value = '2.603'
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("E42").Value = '  +1.50%  '
$ws.Range("D43").Value = '5.628'
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("D44").Value = '0.8171'
$ws.Range("E44").Value = '  -3.96%  '
$ws.Range("D45").Value = '99.80'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").Value = '1.777.42'
$ws.Range("E46").Value = '  -1.87%  '
$ws.Range("D47").Value = '0.0₈114'
$ws.Range("E47").Value = '  -1.51%  '
$ws.Range("D48").Value = '0.4571'
$ws.Range("E48").Value = '  +1.74%  '
$ws.Range("D49").Value = '1.015'
$ws.Range("E49").Value = '  +0.94%  '
$ws.Range("D50").Value = '54.84'
$ws.Range("E50").Value = '  -2.15%  '
$ws.Range("D51").Value = '0.05042'
$ws.Range("E51").Value = '  -1.73%  '
